# Quicker messages update. Pastes rather than typing
#
# The single sample recipient row is replaced with the real pasted
# WhatsApp recipient/group list (rows 2-13).
#
# Values go in first, row by row, left to right (A, B, C, E, F; D is a
# plain paste-count number) so freshly-introduced text lands in the same
# shared-string slots the source workbook ends up with. F6 ('Ms Shan') is
# written dead last because it is the one brand-new string that isn't
# introduced in row order.
#
# Hyperlinks are then wired up in a second pass, in the same left-to-right
# sweep-with-a-backtrack order the relationship ids show up in the saved
# file (B3, B6, B8, B9, B7, B11) -- emails that came in already linked
# keep plain formatting (B3); the rest pick up the Hyperlink cell style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Cell values ----------------------------------------------------

# Row 2
$ws.Range("A2").Value = "Kavish"
$ws.Range("B2").Value = "kvishrock120809@gmail.com"
$ws.Range("C2").Value = "+60125262136"
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = "GFlriDI5wB4FhyBqCygwPB"
$ws.Range("F2").Value = "Kavish & Hehmaa"

# Row 3
$ws.Range("A3").Value = "Hehmaa"
$ws.Range("B3").Value = "kaushi180370@gmail.com"
$ws.Range("C3").Value = "+60125262136"
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = "GFlriDI5wB4FhyBqCygwPB"
$ws.Range("F3").Value = "Kavish & Hehmaa"

# Row 4
$ws.Range("A4").Value = "David"
$ws.Range("B4").Value = "dss.batuncang@gmail.com"
$ws.Range("C4").Value = "+60198188727"
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = "DxsC8kOzmM80ZUKgPGzrIw"
$ws.Range("F4").Value = "all"

# Row 5
$ws.Range("A5").Value = "Samuil-Ashton"
$ws.Range("B5").Value = "samsatu@gmail.com"
$ws.Range("C5").Value = "+60138362370"
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = "DxsC8kOzmM80ZUKgPGzrIw"
$ws.Range("F5").Value = "all"

# Row 6
$ws.Range("A6").Value = "Shan"
$ws.Range("B6").Value = "shanredai@gmail.com"
$ws.Range("C6").Value = "+919740081143"
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = "BXVVhuHWqZkJwOaM1jARrG"

# Row 7
$ws.Range("A7").Value = "Priyakari"
$ws.Range("B7").Value = "ramadevu.priyakari@gmail.com"
$ws.Range("C7").Value = "+919441788160"
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = "JjqgTyDwbxfFUQsHfAuoty"

# Row 8
$ws.Range("A8").Value = "Sudiksha"
$ws.Range("B8").Value = "chsudiksha3@gmail.com"
$ws.Range("C8").Value = "+919880874620"
$ws.Range("D8").Value = 1
$ws.Range("E8").Value = "HEnZ5UBNjh60r2fsAU3Ci4"

# Row 9
$ws.Range("A9").Value = "Satwik"
$ws.Range("B9").Value = "apparasu1966@gmail.com"
$ws.Range("C9").Value = "+919866124793"
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = "B5B4TDTMwnx4NZ7xTd3WBE"

# Row 10
$ws.Range("A10").Value = "Kushal"
$ws.Range("B10").Value = "kingkushalraj22@gmail.com"
$ws.Range("C10").Value = "+601169239511"
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = "IUPuVg78fcA6Sh7CQWNxmf"

# Row 11
$ws.Range("A11").Value = "Radha"
$ws.Range("B11").Value = "radhachaganti112@gmail.com"
$ws.Range("C11").Value = "+919491392460"
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = "FOkOvmWcghgDBrm2VUhlOi"
$ws.Range("F11").Value = "all"

# Row 12
$ws.Range("A12").Value = "Guventhra"
$ws.Range("C12").Value = "+60102701163"
$ws.Range("D12").Value = 1
$ws.Range("E12").Value = "FOkOvmWcghgDBrm2VUhlOi"
$ws.Range("F12").Value = "all"

# Row 13
$ws.Range("A13").Value = "Vamsi Krishna"
$ws.Range("C13").Value = "+918500970197"
$ws.Range("D13").Value = 1
$ws.Range("E13").Value = "DVUkMcwQH0k1dlYfAKHvzQ"
$ws.Range("F13").Value = "all"

# F6 is written last: 'Ms Shan' is the very last brand-new shared string
# in the target workbook, landing after every other new cell.
$ws.Range("F6").Value = "Ms Shan"

# --- 2. Clear the leftover 'Hyperlink' look on the cells that weren't --
#        re-linked (they inherited it, empty, from the old template)   --
$ws.Range("B2").Style = "Normal"
$ws.Range("B4").Style = "Normal"
$ws.Range("B5").Style = "Normal"
$ws.Range("B10").Style = "Normal"

# --- 3. Hyperlinks (second pass, matches saved relationship-id order) --
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:kaushi180370@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B6"), "mailto:shanredai@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B8"), "mailto:chsudiksha3@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B9"), "mailto:apparasu1966@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B7"), "mailto:ramadevu.priyakari@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B11"), "mailto:radhachaganti112@gmail.com") | Out-Null

# B3 stays plain; the rest pick up the Hyperlink cell style.
$ws.Range("B3").Style = "Normal"
$ws.Range("B6").Style = "Hyperlink"
$ws.Range("B7").Style = "Hyperlink"
$ws.Range("B8").Style = "Hyperlink"
$ws.Range("B9").Style = "Hyperlink"
$ws.Range("B11").Style = "Hyperlink"

# Cursor/selection after the paste settles one column past the table.
$ws.Range("G16").Select() | Out-Null
